$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row heights for rows 2-4 change from 290 to max (409.5)
$ws.Rows.Item(2).RowHeight = 409.5
$ws.Rows.Item(3).RowHeight = 409.5
$ws.Rows.Item(4).RowHeight = 409.5

# New big stats query replaces the old StatQuery cell content (same text for all 3 rows)
$statQuery = "MATCH (p:program)<--(s:study)<--(c:case)<--(demo:demographic), (c)<--(diag:diagnosis)`n      WHERE (size([]) = 0 OR s.clinical_study_designation IN [])`n        AND (s.study_disposition = 'Unrestricted')`n        AND s.clinical_study_designation IN ['UBC01']`n  and diag.stage_of_disease in [ 'T3N0M1', 'Not Applicable'] `n        AND (size([]) = 0 OR demo.sex IN [])`n        AND (size([]) = 0 OR demo.neutered_indicator IN [])`n        AND (size([]) = 0 OR diag.disease_term IN [])`n        AND (size([]) = 0 OR diag.primary_disease_site IN [])`n        AND (size([]) = 0 OR diag.stage_of_disease IN [])`n        AND (size([]) = 0 OR diag.best_response IN [])`n    OPTIONAL MATCH (c)-->(co:cohort)`n    OPTIONAL MATCH (f:file)-[*]->(c)`n    OPTIONAL MATCH (f)-->(parent)`n    OPTIONAL MATCH (samp:sample)-->(c)`n    OPTIONAL MATCH (samp)<--(al:aliquot)`n    WITH DISTINCT c AS c, p, s, co, demo, diag, f, parent, samp, al`n      WHERE (size([]) = 0 OR samp.summarized_sample_type IN [])`n        AND (size([]) = 0 OR samp.specific_sample_pathology IN [])`n        AND (size([]) = 0 OR samp.sample_site IN [])`n        AND (size([]) = 0 OR head(labels(parent)) IN [])`n        AND (size([]) = 0 OR f.file_type IN [])`n        AND (size([]) = 0 OR f.file_format IN [])`n    WITH c.case_id AS case_id,`n         s.clinical_study_designation AS study_code,`n         s.clinical_study_type AS study_type,`n         co.cohort_description AS cohort,`n         demo.breed AS breed,`n         diag.disease_term AS diagnosis,`n         diag.stage_of_disease AS stage_of_disease,`n         diag.primary_disease_site AS disease_site,`n         demo.patient_age_at_enrollment AS age,`n         demo.sex AS sex,`n         demo.neutered_indicator AS neutered_status,`n         demo.weight AS weight,`n         diag.best_response AS response_to_treatment,`n         samp.sample_id AS sample_id,`n         f.uuid AS file_id,`n         al`n    RETURN`nCOUNT(DISTINCT file_id) as number_of_files,`nCOUNT(DISTINCT sample_id) as number_of_sample,`nCOUNT(DISTINCT case_id) as number_of_cases,`nCOUNT(DISTINCT study_code) as number_of_study,`nCOUNT(DISTINCT al) as number_of_aliquot`n    "
$ws.Range("C2").Value = $statQuery
$ws.Range("C3").Value = $statQuery
$ws.Range("C4").Value = $statQuery

# The query lines 3-46 (skipping first 2 lines) are also written individually down column C starting row 5
$ws.Range("C5").Value = "        AND (s.study_disposition = 'Unrestricted')"
$ws.Range("C6").Value = "        AND s.clinical_study_designation IN ['UBC01']"
$ws.Range("C7").Value = "  and diag.stage_of_disease in [ 'T3N0M1', 'Not Applicable'] "
$ws.Range("C8").Value = "        AND (size([]) = 0 OR demo.sex IN [])"
$ws.Range("C9").Value = "        AND (size([]) = 0 OR demo.neutered_indicator IN [])"
$ws.Range("C10").Value = "        AND (size([]) = 0 OR diag.disease_term IN [])"
$ws.Range("C11").Value = "        AND (size([]) = 0 OR diag.primary_disease_site IN [])"
$ws.Range("C12").Value = "        AND (size([]) = 0 OR diag.stage_of_disease IN [])"
$ws.Range("C13").Value = "        AND (size([]) = 0 OR diag.best_response IN [])"
$ws.Range("C14").Value = "    OPTIONAL MATCH (c)-->(co:cohort)"
$ws.Range("C15").Value = "    OPTIONAL MATCH (f:file)-[*]->(c)"
$ws.Range("C16").Value = "    OPTIONAL MATCH (f)-->(parent)"
$ws.Range("C17").Value = "    OPTIONAL MATCH (samp:sample)-->(c)"
$ws.Range("C18").Value = "    OPTIONAL MATCH (samp)<--(al:aliquot)"
$ws.Range("C19").Value = "    WITH DISTINCT c AS c, p, s, co, demo, diag, f, parent, samp, al"
$ws.Range("C20").Value = "      WHERE (size([]) = 0 OR samp.summarized_sample_type IN [])"
$ws.Range("C21").Value = "        AND (size([]) = 0 OR samp.specific_sample_pathology IN [])"
$ws.Range("C22").Value = "        AND (size([]) = 0 OR samp.sample_site IN [])"
$ws.Range("C23").Value = "        AND (size([]) = 0 OR head(labels(parent)) IN [])"
$ws.Range("C24").Value = "        AND (size([]) = 0 OR f.file_type IN [])"
$ws.Range("C25").Value = "        AND (size([]) = 0 OR f.file_format IN [])"
$ws.Range("C26").Value = "    WITH c.case_id AS case_id,"
$ws.Range("C27").Value = "         s.clinical_study_designation AS study_code,"
$ws.Range("C28").Value = "         s.clinical_study_type AS study_type,"
$ws.Range("C29").Value = "         co.cohort_description AS cohort,"
$ws.Range("C30").Value = "         demo.breed AS breed,"
$ws.Range("C31").Value = "         diag.disease_term AS diagnosis,"
$ws.Range("C32").Value = "         diag.stage_of_disease AS stage_of_disease,"
$ws.Range("C33").Value = "         diag.primary_disease_site AS disease_site,"
$ws.Range("C34").Value = "         demo.patient_age_at_enrollment AS age,"
$ws.Range("C35").Value = "         demo.sex AS sex,"
$ws.Range("C36").Value = "         demo.neutered_indicator AS neutered_status,"
$ws.Range("C37").Value = "         demo.weight AS weight,"
$ws.Range("C38").Value = "         diag.best_response AS response_to_treatment,"
$ws.Range("C39").Value = "         samp.sample_id AS sample_id,"
$ws.Range("C40").Value = "         f.uuid AS file_id,"
$ws.Range("C41").Value = "         al"
$ws.Range("C42").Value = "    RETURN"
$ws.Range("C43").Value = "COUNT(DISTINCT file_id) as number_of_files,"
$ws.Range("C44").Value = "COUNT(DISTINCT sample_id) as number_of_sample,"
$ws.Range("C45").Value = "COUNT(DISTINCT case_id) as number_of_cases,"
$ws.Range("C46").Value = "COUNT(DISTINCT study_code) as number_of_study,"
$ws.Range("C47").Value = "COUNT(DISTINCT al) as number_of_aliquot"
$ws.Range("C48").Value = "    "

# Update active selection to C2
$ws.Range("C2").Select()
